# Update peakValue (G), RMS (H), tau (I), and AUC (J) columns for rows 2-19
# with recalculated values per "uniformity test, tester fix"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 50.15000000000001
$ws.Range("H2").Value = 12.84536111111111
$ws.Range("I2").Value = 153.9534722222222
$ws.Range("J2").Value = 5027.915722222222
$ws.Range("G3").Value = 81.98425925925926
$ws.Range("H3").Value = 22.06077777777778
$ws.Range("I3").Value = 156.9860740740741
$ws.Range("J3").Value = 9333.02887037037
$ws.Range("G4").Value = 66.86018518518519
$ws.Range("H4").Value = 19.27598148148148
$ws.Range("I4").Value = 199.1854907407407
$ws.Range("J4").Value = 8785.400435185185
$ws.Range("G5").Value = 67.20555555555556
$ws.Range("H5").Value = 18.26519444444445
$ws.Range("I5").Value = 162.8410555555556
$ws.Range("J5").Value = 7508.073027777778
$ws.Range("G6").Value = 114.0638888888889
$ws.Range("H6").Value = 32.93819444444445
$ws.Range("I6").Value = 176.8305740740741
$ws.Range("J6").Value = 14621.79834259259
$ws.Range("G7").Value = 99.12037037037038
$ws.Range("H7").Value = 26.91130555555556
$ws.Range("I7").Value = 145.1042222222222
$ws.Range("J7").Value = 10822.44786111111
$ws.Range("G8").Value = 123.8351851851852
$ws.Range("H8").Value = 33.65185185185185
$ws.Range("I8").Value = 150.8492407407408
$ws.Range("J8").Value = 13712.78072222222
$ws.Range("G9").Value = 251.6518518518519
$ws.Range("H9").Value = 68.40933333333334
$ws.Range("I9").Value = 164.5025370370371
$ws.Range("J9").Value = 27733.00499074074
$ws.Range("G10").Value = 139.937962962963
$ws.Range("H10").Value = 33.95793518518519
$ws.Range("I10").Value = 97.69907407407408
$ws.Range("J10").Value = 11767.41068518518
$ws.Range("G11").Value = 143.3435185185185
$ws.Range("H11").Value = 38.61753703703704
$ws.Range("I11").Value = 137.7861018518519
$ws.Range("J11").Value = 14239.41121296296
$ws.Range("G12").Value = 216.7611111111111
$ws.Range("H12").Value = 49.80016666666667
$ws.Range("I12").Value = 69.07636111111111
$ws.Range("J12").Value = 15780.26272222222
$ws.Range("G13").Value = 184.6898148148148
$ws.Range("H13").Value = 50.48237037037037
$ws.Range("I13").Value = 161.3433888888889
$ws.Range("J13").Value = 20323.01669444445
$ws.Range("G14").Value = 105.9212962962963
$ws.Range("H14").Value = 26.42299074074074
$ws.Range("I14").Value = 116.863537037037
$ws.Range("J14").Value = 9238.739981481482
$ws.Range("G15").Value = 99.17129629629629
$ws.Range("H15").Value = 23.67728703703704
$ws.Range("I15").Value = 117.0621759259259
$ws.Range("J15").Value = 8837.139092592592
$ws.Range("G16").Value = 151.0166666666667
$ws.Range("H16").Value = 30.1268425925926
$ws.Range("I16").Value = 60.95931481481481
$ws.Range("J16").Value = 9303.403685185185
$ws.Range("G17").Value = 72.50555555555555
$ws.Range("H17").Value = 18.51413888888889
$ws.Range("I17").Value = 153.5498055555555
$ws.Range("J17").Value = 7226.979
$ws.Range("G18").Value = 122.9666666666667
$ws.Range("H18").Value = 34.91819444444445
$ws.Range("I18").Value = 144.6688333333333
$ws.Range("J18").Value = 14414.89663888889
$ws.Range("G19").Value = 119.6296296296296
$ws.Range("H19").Value = 29.74342592592593
$ws.Range("I19").Value = 121.2216944444444
$ws.Range("J19").Value = 10485.57519444444
